$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Age column (B) holds numeric-looking text in this sheet (e.g. "57", "54"),
# so force text formatting before assignment to avoid auto-conversion to a number.
$ws.Range("B10:B11").NumberFormat = "@"

$ws.Range("A10").Value = "Hero Pulsur"
$ws.Range("B10").Value = "100"
$ws.Range("C10").Value = "Other"
$ws.Range("D10").Value = "08-01-2025 20:59:27"
$ws.Range("E10").Value = "working"
$ws.Range("F10").Value = "working"
$ws.Range("G10").Value = "working"
$ws.Range("H10").Value = "working"

$ws.Range("A11").Value = "Swapnanil Bala"
$ws.Range("B11").Value = "27"
$ws.Range("C11").Value = "Male"
$ws.Range("D11").Value = "08-01-2025 21:00:59"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "not"
$ws.Range("G11").Value = "working"
$ws.Range("H11").Value = "yet"
